# Dummy server service implementation.
#
# - Mark "Create server service" (row 7) as Done.
# - Rename the old "Create build strategy" task (row 9) to
#   "Create project sturcture and build strategy".
# - Add a new row 10 task "Integrate server service with database" (pending).
# - Column B auto-fits to the new, longer text.
# - Move the active selection to F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("Create server service") moves from pending -> done:
# copy the formatting used by the other "Done" rows (e.g. row 6) onto row 7,
# then fill in the Status value.
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C7").Value = "Done"

# Row 9: task text changes (still "Done").
$ws.Range("B9").Value = "Create project sturcture and build strategy"

# New row 10: copy the pending-row formatting (row 8) down one row, then set values.
$ws.Range("A8:C8").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Integrate server service with database"

# The Status column's font is normally size 9; the new empty status cell on
# row 10 uses the regular size-11 variant of the same pending color/fill.
$ws.Range("C10").Font.Size = 11

$ws.Columns.Item(2).AutoFit()

$ws.Range("F9").Select()
